$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.850.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.55%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.619.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.17%  '

# Row 4
$ws.Range("E4").Value = '  -1.04%  '

# Row 5
$ws.Range("E5").Value = '  -2.03%  '

# Row 6
$ws.Range("E6").Value = '  -1.79%  '

# Row 7
$ws.Range("E7").Value = '  -1.08%  '

# Row 8
$ws.Range("E8").Value = '  -2.16%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0616'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.65%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.32%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.34%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.844.12'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.17%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.613.75'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.41%  '

# Row 14
$ws.Range("E14").Value = '  -3.16%  '

# Row 15
$ws.Range("E15").Value = '  -3.49%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.866.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.56%  '

# Row 17
$ws.Range("E17").Value = '  -3.34%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0736'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.67%  '

# Row 19
$ws.Range("E19").Value = '  -1.06%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.68%  '

# Row 22
$ws.Range("E22").Value = '  -2.97%  '

# Row 23
$ws.Range("E23").Value = '  -2.82%  '

# Row 24
$ws.Range("E24").Value = '  +2.22%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.70'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.98%  '

# Row 26
$ws.Range("E26").Value = '  -1.34%  '

# Row 27
$ws.Range("E27").Value = '  -3.54%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.45%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.30%  '

# Row 30
$ws.Range("E30").Value = '  -2.04%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0477'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.61%  '

# Row 32
$ws.Range("E32").Value = '  -4.39%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.09'
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = '  -2.67%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.51%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.124.52'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.26%  '

# Row 37
$ws.Range("E37").Value = '  -7.35%  '

# Row 39
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.511'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.20%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0153'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.57%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.68%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.755.13'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.05%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.748'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.89%  '

# Row 44
$ws.Range("E44").Value = '  -5.91%  '

# Row 45
$ws.Range("E45").Value = '  -2.37%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.56%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.59%  '

# Row 48
$ws.Range("E48").Value = '  -0.75%  '

# Row 49
$ws.Range("E49").Value = '  -1.85%  '

# Row 50
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.98%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.85%  '
